# Update the "Student names" credits textbox on slide 1:
#  - split the existing "...Utkarsh khuspare" run so "khuspare" is its own run
#  - append ", Praveg Chikte," to that paragraph (as two more runs)
#  - add a new paragraph "\t\tAryan Raut"
#  - grow the textbox to fit the extra line

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 1")
$tf  = $shp.TextFrame
$tr  = $tf.TextRange

# Paragraph 2 currently reads "Student names :Utkarsh khuspare"
$para2 = $tr.Paragraphs(2, 1)

# Split "khuspare" off into its own run (re-assigning identical text to a
# sub-range forces the run to be split at that boundary).
$splitAt = $para2.Text.IndexOf("khuspare")
$nameRange = $tr.Characters($para2.Start + $splitAt, $para2.Length - $splitAt)
$nameRange.Text = "khuspare"

# Append the rest of the first line as two more runs on the same paragraph.
$para2 = $tr.Paragraphs(2, 1)
$r1 = $para2.InsertAfter(", Praveg ")
$r2 = $r1.InsertAfter("Chikte,")

# Start a new paragraph (carriage return) with the second line of names.
[void]$r2.InsertAfter("`r`t`tAryan Raut")

# The extra line needs more room in the (auto-fit) text box.
$shp.Height = [double]75.1417360
